$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thông tin sản phẩm")
$ws.Activate()

# Drop the stray empty trailing row left over below the data.
$ws.Rows("8").Delete()

# Insert a new column D ("Giá niêm yết" / listed price), shifting the
# existing Đơn vị..Địa chỉ kho columns (D:J) right to (E:K).
$ws.Columns("D").Insert()

# Header
$ws.Range("D1").Value = "Giá niêm yết"
$ws.Range("I1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Data rows - new listed-price values
$ws.Range("D2").Value = 50000
$ws.Range("D3").Value = 30000
$ws.Range("D4").Value = 30000
$ws.Range("D5").Value = 40000
$ws.Range("I2:I5").Copy()
$ws.Range("D2:D5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("K2").Select()
$ws.Range("A1:K2").Select()
